$wb = $excel.ActiveWorkbook

# "zh-cn" worksheet: update Correspond Handoff/Handback DateTime for row 7
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D7").Value = "2016-03-04 10:40:28"
$wsZh.Range("G7").Value = "2016-03-04 10:41:56"

# "de-de" worksheet: update Correspond Handoff/Handback DateTime for row 7
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D7").Value = "2016-03-04 10:40:46"
$wsDe.Range("G7").Value = "2016-03-04 10:42:22"
